# Fix expected xlsx for exportServiceIT
#
# This reproduces the data changes from the commit:
#  - "packages" sheet: a handful of cells that only held the placeholder
#    empty string are cleared out (C2, D2, C3, C4, D4)
#  - "attributes" sheet: column O placeholder empty-string cells are
#    cleared out for (almost) every data row (O10 keeps its real value)
#  - "pack_test1" sheet: the date_attr (F) and datetime_attr (H) columns
#    were exported as formatted text instead of raw Excel serial numbers

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "packages": drop cells that only contained the empty placeholder
# ---------------------------------------------------------------------
$wsPackages = $wb.Worksheets.Item("packages")
$wsPackages.Range("C2").ClearContents()
$wsPackages.Range("D2").ClearContents()
$wsPackages.Range("C3").ClearContents()
$wsPackages.Range("C4").ClearContents()
$wsPackages.Range("D4").ClearContents()

# ---------------------------------------------------------------------
# Sheet "attributes": drop column O placeholder cells for every row
# except row 10 (which legitimately holds "option1,option2")
# ---------------------------------------------------------------------
$wsAttributes = $wb.Worksheets.Item("attributes")
$oRows = 2..25 | Where-Object { $_ -ne 10 }
foreach ($r in $oRows) {
    $wsAttributes.Range("O$r").ClearContents()
}

# ---------------------------------------------------------------------
# Sheet "pack_test1": date_attr / datetime_attr columns become text
# ---------------------------------------------------------------------
$wsPackTest1 = $wb.Worksheets.Item("pack_test1")

$dateCells = @("F2", "F3", "F4", "F5")
foreach ($cell in $dateCells) {
    $rng = $wsPackTest1.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = "2012-12-21"
    $rng.ClearFormats()
}

$dateTimeCells = @("H2", "H3", "H4", "H5")
foreach ($cell in $dateTimeCells) {
    $rng = $wsPackTest1.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = "1985-08-12T08:12:13+0200"
    $rng.ClearFormats()
}
